# Add a new "Cart Testdata" worksheet after the existing "Testdata" sheet,
# populate it with the new test-case row, size its columns, and move the
# active tab / selections to match the authored change.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet is inserted right after "Testdata" -> becomes the 2nd sheet
# (sheetId 6, matching workbook.xml) and automatically becomes the active tab.
$cartSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$cartSheet.Name = "Cart Testdata"

# New test-data row (TC Name / Pass-Fail-Skip header + one result row).
$cartSheet.Range("A1").Value = "TC Name"
$cartSheet.Range("B1").Value = "Pass/Fail/Skip"
$cartSheet.Range("A2").Value = "Validate Cart Message"
$cartSheet.Range("B2").Value = "Pass"

# Column widths (best-fit sized columns on the new sheet).
$cartSheet.Columns.Item(1).ColumnWidth = 19.17
$cartSheet.Columns.Item(2).ColumnWidth = 12.67

# Selection bookkeeping: old sheet's cursor moves off its old A5 selection,
# new sheet becomes selected/active with its own cursor at E7.
$ws1.Range("A2").Select() | Out-Null
$cartSheet.Range("E7").Select() | Out-Null
